$d = $word.ActiveDocument

$d.Content.Find.Execute(" (2019)", $true, $false, $false, $false, $false,
                         $true, 1, $false, " (2019) ", 2)
$d.Content.Find.Execute(" (2017)", $true, $false, $false, $false, $false,
                         $true, 1, $false, " (2017) ", 2)
$d.Content.Find.Execute(" (2006)", $true, $false, $false, $false, $false,
                         $true, 1, $false, " (2006) ", 2)
$d.Content.Find.Execute(" (2004)", $true, $false, $false, $false, $false,
                         $true, 1, $false, " (2004) ", 2)
